$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Cells that change FROM a text placeholder TO a numeric value, or vice versa between numeric styles ---
# First: copy-paste exact formatting from stable donor cells, to replicate target cell style exactly.

# Text-style target cells (copy text+style from donor; no value set needed afterward)
$ws.Range("D14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("D14").Copy($ws.Range("F15"))
$ws.Range("D14").Copy($ws.Range("C26"))
$ws.Range("D14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("D14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("D14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))

# Numeric-style target cells (copy style from donor, then set the new numeric value)
$ws.Range("J14").Copy($ws.Range("C14"))
$ws.Range("C14").Value = 1
$ws.Range("J14").Copy($ws.Range("F14"))
$ws.Range("F14").Value = 1
$ws.Range("J14").Copy($ws.Range("I14"))
$ws.Range("I14").Value = 1
$ws.Range("J14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 2
$ws.Range("K15").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("J14").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 2
$ws.Range("K15").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100
$ws.Range("J14").Copy($ws.Range("J30"))
$ws.Range("J30").Value = 2
$ws.Range("K15").Copy($ws.Range("K30"))
$ws.Range("K30").Value = -100

# --- Value-only cell updates (style unchanged) ---
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -85.714285714285
$ws.Range("H15").Value = -100
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -8.695652173913
$ws.Range("I16").Value = 42
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = -17.647058823529
$ws.Range("L16").Value = -12.5
$ws.Range("M16").Value = -16
$ws.Range("N16").Value = -87.037037037037
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 17.777777777777
$ws.Range("I17").Value = 88
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = -2.222222222222
$ws.Range("L17").Value = 54.385964912280
$ws.Range("M17").Value = 83.333333333333
$ws.Range("N17").Value = -30.708661417322
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -17.241379310344
$ws.Range("M18").Value = -31.428571428571
$ws.Range("N18").Value = -92.307692307692
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -27.272727272727
$ws.Range("I19").Value = 68
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = -15
$ws.Range("L19").Value = 23.636363636363
$ws.Range("M19").Value = 119.354838709677
$ws.Range("N19").Value = -45.6
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 18
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 39
$ws.Range("J20").Value = 43
$ws.Range("K20").Value = -9.302325581395
$ws.Range("L20").Value = 200
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -72.340425531914
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 137
$ws.Range("G21").Value = 156
$ws.Range("H21").Value = -12.179487179487
$ws.Range("I21").Value = 264
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = -12
$ws.Range("L21").Value = 29.411764705882
$ws.Range("M21").Value = 48.314606741573
$ws.Range("N21").Value = -74.664107485604
$ws.Range("L22").Value = -25
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -58.064516129032
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -27.956989247311
$ws.Range("I24").Value = 117
$ws.Range("J24").Value = 164
$ws.Range("K24").Value = -28.658536585365
$ws.Range("L24").Value = 4.464285714285
$ws.Range("M24").Value = -10
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 99
$ws.Range("G25").Value = 85
$ws.Range("H25").Value = 16.470588235294
$ws.Range("I25").Value = 167
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 39.166666666666
$ws.Range("L25").Value = 57.547169811320
$ws.Range("M25").Value = 18.439716312056
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = 350
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 11
$ws.Range("J27").Value = 11
$ws.Range("L27").Value = 37.5
$ws.Range("C28").Value = 1
$ws.Range("I28").Value = 5
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = -37.5
$ws.Range("N28").Value = -76.190476190476
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = -20
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -78.947368421052
